$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "Computer Software, Biochemistry"
$ws.Range("F2").Select()
